# "Data source corrected and updated"
# Columns J and K (rows 1-51) held a mis-sourced set of values (header row
# used text labels "r"/"s" via shared strings, the data rows held 1 / 0.6).
# The corrected data source uses numeric 0.6 / 0.3 uniformly for every row,
# including the former header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite J1:J51 and K1:K51 with the corrected numeric values. Writing a
# scalar to a multi-cell range fills every cell in that range, so this also
# replaces the old text labels in J1/K1 with plain numbers (dropping the
# shared-string entries) and rewrites every data row from 1/0.6 to 0.6/0.3.
$ws.Range("J1:J51").Value = 0.6
$ws.Range("K1:K51").Value = 0.3

# Restore the on-screen selection to match where the corrected data now
# draws attention (K column) and bring the zoom back to a plain 100%.
$ws.Range("K1:K51").Select()
$win = $excel.ActiveWindow
$win.Zoom = 100
